$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'26.496.66"
$ws.Range("E2").Value = "  -0.24%  "

# Row 3
$ws.Range("D3").Value = "'1.840.22"
$ws.Range("E3").Value = "  -0.34%  "

# Row 4
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").Value = "'261.13"
$ws.Range("E5").Value = "  -1.16%  "

# Row 6
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.13%  "

# Row 7
$ws.Range("D7").Value = "'0.5335"
$ws.Range("E7").Value = "  +2.36%  "

# Row 8
$ws.Range("D8").Value = "'0.3019"
$ws.Range("E8").Value = "  -6.74%  "

# Row 9
$ws.Range("D9").Value = "'0.06880"
$ws.Range("E9").Value = "  +1.09%  "

# Row 10
$ws.Range("D10").Value = "'17.75"
$ws.Range("E10").Value = "  -5.48%  "

# Row 11
$ws.Range("D11").Value = "'1.850.87"
$ws.Range("E11").Value = "  +0.48%  "

# Row 12
$ws.Range("D12").Value = "'0.7361"
$ws.Range("E12").Value = "  -5.41%  "

# Row 13
$ws.Range("D13").Value = "'0.07469"
$ws.Range("E13").Value = "  -3.88%  "

# Row 14
$ws.Range("D14").Value = "'89.49"
$ws.Range("E14").Value = "  +1.09%  "

# Row 15
$ws.Range("D15").Value = "'4.978"
$ws.Range("E15").Value = "  -1.03%  "

# Row 16
$ws.Range("E16").Value = "  +0.17%  "

# Row 17
$ws.Range("D17").Value = "'13.95"
$ws.Range("E17").Value = "  -0.31%  "

# Row 18
$ws.Range("E18").Value = "  +0.11%  "

# Row 19
$ws.Range("D19").Value = "'0.000007918"
$ws.Range("E19").Value = "  -0.55%  "

# Row 20
$ws.Range("D20").Value = "'26.527.58"
$ws.Range("E20").Value = "  -0.32%  "

# Row 21
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'4.594"
$ws.Range("E21").Value = "  -0.83%  "

# Row 22
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").Value = "'5.979"
$ws.Range("E22").Value = "  -0.59%  "

# Row 23
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "'9.286"
$ws.Range("E23").Value = "  -1.86%  "

# Row 24
$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").Value = "'142.91"
$ws.Range("E24").Value = "  -0.35%  "

# Row 25
$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D25").Value = "'2.220"
$ws.Range("E25").Value = "  +2.02%  "

# Row 26
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'1.688"
$ws.Range("E26").Value = "  +0.54%  "

# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'16.94"
$ws.Range("E27").Value = "  -0.38%  "

# Row 28
$ws.Range("B28").Value = "BitcoinCash"
$ws.Range("C28").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D28").Value = "'110.62"
$ws.Range("E28").Value = "  -1.01%  "

# Row 29
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'4.260"
$ws.Range("E29").Value = "  +1.60%  "

# Row 30
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "'0.08796"
$ws.Range("E30").Value = "  +0.56%  "

# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.036"
$ws.Range("E31").Value = "  -1.93%  "

# Row 32
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.04792"
$ws.Range("E32").Value = "  -1.04%  "

# Row 33
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.7331"
$ws.Range("E33").Value = "  +1.53%  "

# Row 34
$ws.Range("D34").Value = "'2.919"
$ws.Range("E34").Value = "  +2.07%  "

# Row 35
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.130"
$ws.Range("E35").Value = "  -0.14%  "

# Row 36
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").Value = "'3.099"
$ws.Range("E36").Value = "  -0.20%  "

# Row 37
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "'2.291"
$ws.Range("E37").Value = "  +3.24%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01714"
$ws.Range("E38").Value = "  -4.40%  "

# Row 39
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.4735"
$ws.Range("E39").Value = "  -2.60%  "

# Row 40
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'0.9026"
$ws.Range("E40").Value = "  +1.00%  "

# Row 41
$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").Value = "'107.69"
$ws.Range("E41").Value = "  -3.19%  "

# Row 42
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'5.869"
$ws.Range("E42").Value = "  -2.78%  "

# Row 43
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "'1.001"
$ws.Range("E43").Value = "  +0.13%  "

# Row 44
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "'7.412"
$ws.Range("E44").Value = "  -2.73%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'9.034"
$ws.Range("E45").Value = "  -0.40%  "

# Row 46
$ws.Range("D46").Value = "'0.4095"
$ws.Range("E46").Value = "  -2.78%  "

# Row 47
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "'0.1235"
$ws.Range("E47").Value = "  -0.33%  "

# Row 48
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "'34.92"
$ws.Range("E48").Value = "  -0.35%  "

# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.05798"
$ws.Range("E49").Value = "  -1.56%  "

# Row 50
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "'0.8953"
$ws.Range("E50").Value = "  +0.86%  "

# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'60.20"
$ws.Range("E51").Value = "  +0.43%  "

